# "ready benchmark for Lockfree queues" — append the 19 марта diary entry
# to the Лист1 work-log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 24 was the next blank row in the log; fill in the date and the
# description of the work done.
$ws.Cells.Item(24, 1).Value = "19 марта"
$ws.Cells.Item(24, 2).Value = "Правка теоритической части текста в черновике вкр. Создание проекта для проведения бенчмарков различных библиотек, а также написание самих бенчмарков"

# Column B auto-sizes ("best fit") to the widest entry; recompute it now
# that a wider (mostly Cyrillic) line has been added.
$ws.Columns.Item(2).AutoFit() | Out-Null

# The saved selection moves on to the next empty row, column A.
$ws.Range("A25").Select()
